$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing existing rows 10-31 down to 11-32
$ws.Rows.Item(10).Insert()

# Fill the new row 10 with values (copy of the constant columns from the
# surrounding rows, plus the new data point values)
$ws.Cells.Item(10, 1).Value = 2
$ws.Cells.Item(10, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(10, 3).Value = "Coquimbo"
$ws.Cells.Item(10, 4).Value = 44497
$ws.Cells.Item(10, 4).Style = $ws.Cells.Item(11, 4).Style
$ws.Cells.Item(10, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat
$ws.Cells.Item(10, 5).Value = 4
$ws.Cells.Item(10, 6).Value = 100112032
$ws.Cells.Item(10, 7).Value = "Zapallo italiano"
$ws.Cells.Item(10, 8).Value = "Sin especificar"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 800
$ws.Cells.Item(10, 11).Value = 7500
$ws.Cells.Item(10, 12).Value = 8000
$ws.Cells.Item(10, 13).Value = 7750
$ws.Cells.Item(10, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 129
$ws.Cells.Item(10, 17).Value = 60
$ws.Cells.Item(10, 18).Value = "Hortaliza"
